$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (rows 4 and 5), pushing the
# existing data down.
$ws.Rows("4:5").Insert()

$ws.Range("A4").Value = "好想爱这个世界啊"
$ws.Range("B4").Value = "华晨宇"
$ws.Range("C4").Value = "eSvnAyHFoyo"

$ws.Range("A5").Value = "我管你"
$ws.Range("B5").Value = "华晨宇"
$ws.Range("C5").Value = "wmRkAWPuvCg"

# Insert one new row after "回不去的夏天" (now at row 10), pushing the rest
# of the data down.
$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = "尘星"
$ws.Range("B11").Value = "夏日入侵企画"
$ws.Range("C11").Value = "jS0rDCTyg-E"

# Update selection to match the saved workbook view.
$ws.Range("E12").Select()
